$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "62.654.31"
Set-TextValue "E2" "  -3.87%  "
Set-TextValue "D3" "3.043.19"
Set-TextValue "E3" "  -3.25%  "
Set-TextValue "E4" "  +0.16%  "
Set-TextValue "D5" "543.15"
Set-TextValue "E5" "  -4.45%  "
Set-TextValue "D6" "133.55"
Set-TextValue "E6" "  -10.72%  "
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  +0.09%  "
Set-TextValue "D8" "3.039.74"
Set-TextValue "E8" "  -3.09%  "
Set-TextValue "D9" "0.488"
Set-TextValue "E9" "  -2.90%  "
Set-TextValue "B10" "Dogecoin"
Set-TextValue "C10" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue "D10" "0.154"
Set-TextValue "E10" "  -4.17%  "
Set-TextValue "B11" "Toncoin"
Set-TextValue "C11" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D11" "6.34"
Set-TextValue "E11" "  -11.69%  "
Set-TextValue "D12" "0.456"
Set-TextValue "E12" "  -2.62%  "
Set-TextValue "D13" "34.52"
Set-TextValue "E13" "  -5.07%  "
Set-TextValue "E14" "  -5.54%  "
Set-TextValue "D15" "3.542.25"
Set-TextValue "E15" "  -2.96%  "
Set-TextValue "D16" "62.760.55"
Set-TextValue "E16" "  -3.68%  "
Set-TextValue "D18" "3.042.75"
Set-TextValue "E18" "  -3.21%  "
Set-TextValue "E19" "  -3.43%  "
Set-TextValue "D20" "478.68"
Set-TextValue "E20" "  -12.05%  "
Set-TextValue "D21" "13.30"
Set-TextValue "E21" "  -4.91%  "
Set-TextValue "E22" "  -2.63%  "
Set-TextValue "D23" "6.98"
Set-TextValue "E23" "  -6.83%  "
Set-TextValue "D24" "76.94"
Set-TextValue "E24" "  -3.05%  "
Set-TextValue "D25" "12.11"
Set-TextValue "E25" "  -6.91%  "
Set-TextValue "D26" "0.998"
Set-TextValue "E26" "  -0.33%  "
Set-TextValue "E27" "  -4.71%  "
Set-TextValue "E28" "  -9.01%  "
Set-TextValue "E29" "  -0.13%  "
Set-TextValue "E30" "  -11.37%  "
Set-TextValue "D31" "26.03"
Set-TextValue "E31" "  -1.65%  "
Set-TextValue "E32" "  -3.31%  "
Set-TextValue "B33" "Stacks"
Set-TextValue "C33" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D33" "2.47"
Set-TextValue "E33" "  -7.78%  "
Set-TextValue "B34" "OKB"
Set-TextValue "C34" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D34" "59.13"
Set-TextValue "E34" "  +11.79%  "
Set-TextValue "D35" "506.55"
Set-TextValue "E35" "  -8.65%  "
Set-TextValue "D36" "5.91"
Set-TextValue "E36" "  -4.03%  "
Set-TextValue "D37" "5.04"
Set-TextValue "E37" "  -8.32%  "
Set-TextValue "E38" "  -12.47%  "
Set-TextValue "D39" "3.068.32"
Set-TextValue "E39" "  -0.12%  "
Set-TextValue "D40" "0.0782"
Set-TextValue "E40" "  -5.48%  "
Set-TextValue "E41" "  -4.82%  "
Set-TextValue "D42" "7.98"
Set-TextValue "E42" "  -4.46%  "
Set-TextValue "D43" "2.56"
Set-TextValue "E43" "  -12.34%  "
Set-TextValue "E44" "  -4.68%  "
Set-TextValue "E45" "  +0.00%  "
Set-TextValue "E46" "  -9.79%  "
Set-TextValue "D47" "24.34"
Set-TextValue "E47" "  -3.83%  "
Set-TextValue "D48" "119.25"
Set-TextValue "E48" "  -0.69%  "
Set-TextValue "E49" "  -3.38%  "
Set-TextValue "D50" "0.0₃0489"
Set-TextValue "E50" "  -7.95%  "
Set-TextValue "E51" "  +58.81%  "
